# Add a new row for the "hystrix-test" service below the existing entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "hystrix-test"
$ws.Range("B4").Value = 8004

# Move the active selection to A5, just below the newly added row.
$ws.Range("A5").Select()
